$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 43, pushing existing rows 43-57 down to 44-58.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record.
$ws.Range("A43").Value = 11
$ws.Range("B43").Value = "Vega Monumental Concepción"
$ws.Range("C43").Value = "Bíobío"
$ws.Range("D43").Value = "2022-08-03"
$ws.Range("E43").Value = 8
$ws.Range("F43").Value = 100112013
$ws.Range("G43").Value = "Alcachofa"
$ws.Range("H43").Value = "Española"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 150
$ws.Range("K43").Value = 17000
$ws.Range("L43").Value = 18000
$ws.Range("M43").Value = 17467
$ws.Range("N43").Value = "$/caja 30 unidades"
$ws.Range("O43").Value = "Provincia de Limarí"
$ws.Range("P43").Value = 582
$ws.Range("Q43").Value = 30
$ws.Range("R43").Value = "Hortaliza"
